$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 353.73
$ws.Range("D2").Value = 3.76
$ws.Range("E2").Value = 301.89
$ws.Range("F2").Value = 3.75

# Row 3
$ws.Range("C3").Value = 455.93
$ws.Range("D3").Value = 5.96
$ws.Range("F3").Value = 2.48

# Row 5
$ws.Range("E5").Value = 532.22
$ws.Range("F5").Value = 3.47

# Row 7
$ws.Range("E7").Value = 941.36

# Row 8
$ws.Range("E8").Value = 526.5
$ws.Range("F8").Value = 6.58

# Row 9
$ws.Range("E9").Value = 785.46
$ws.Range("F9").Value = 6.35

# Row 10
$ws.Range("E10").Value = 1383

# Row 11
$ws.Range("E11").Value = 147.34
$ws.Range("F11").Value = 2
